# Added checkWinType in basicWins
# The basic-wins reference data (symbol + reel weight table) is reordered
# to match the new win-type check ordering used by checkWinType.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 21
$lastCol = 6

# Snapshot the current A2:F21 block (symbol, reel1..reel5) before shuffling it.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $rowVals
}

# New row order: new row number -> old row number the data comes from.
$newOrder = @{
    2  = 13
    3  = 8
    4  = 10
    5  = 3
    6  = 7
    7  = 5
    8  = 14
    9  = 11
    10 = 9
    11 = 4
    12 = 15
    13 = 2
    14 = 6
    15 = 12
    16 = 19
    17 = 16
    18 = 20
    19 = 18
    20 = 21
    21 = 17
}

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $srcRow = $newOrder[$r]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c - 1]
    }
}
